$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the "Price" column values are numeric-looking strings (e.g. "620.62")
# that the source workbook stores as literal text (to preserve the site's own
# thousands-dot formatting, e.g. "69.013.06"). Force those specific cells to Text
# format before writing so Excel does not auto-convert them to numbers, then clear
# the temporary formatting again so the cells end up with their original (default)
# style, only the content having changed.
$textCells = @("D5", "D6", "D13", "D14", "D18", "D21", "D22", "D25", "D26", "D27", "D29", "D33", "D34", "D35", "D38", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values, row by row (matches the refreshed cryptos snapshot).
$ws.Range("D2").Value = '68.999.25'
$ws.Range("E2").Value = '  +1.32%  '

$ws.Range("D3").Value = '3.765.43'
$ws.Range("E3").Value = '  -0.60%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '620.62'
$ws.Range("E5").Value = '  +3.61%  '

$ws.Range("D6").Value = '164.03'

$ws.Range("D7").Value = '3.763.13'
$ws.Range("E7").Value = '  -0.62%  '

$ws.Range("E8").Value = '  -0.07%  '

$ws.Range("E9").Value = '  +0.65%  '

$ws.Range("E10").Value = '  +1.32%  '

$ws.Range("E11").Value = '  +0.40%  '

$ws.Range("E12").Value = '  +0.79%  '

$ws.Range("D13").Value = '0.0000246'
$ws.Range("E13").Value = '  -0.50%  '

$ws.Range("D14").Value = '35.29'

$ws.Range("D15").Value = '4.403.16'
$ws.Range("E15").Value = '  -0.47%  '

$ws.Range("D16").Value = '3.757.49'
$ws.Range("E16").Value = '  -1.11%  '

$ws.Range("D17").Value = '69.013.06'
$ws.Range("E17").Value = '  +1.37%  '

$ws.Range("D18").Value = '17.66'
$ws.Range("E18").Value = '  -3.38%  '

$ws.Range("E19").Value = '  +0.57%  '

$ws.Range("E20").Value = '  -1.37%  '

$ws.Range("D21").Value = '467.44'
$ws.Range("E21").Value = '  +1.60%  '

$ws.Range("D22").Value = '9.54'
$ws.Range("E22").Value = '  -1.14%  '

$ws.Range("E24").Value = '  +2.52%  '

$ws.Range("D25").Value = '83.00'
$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("D26").Value = '11.98'
$ws.Range("E26").Value = '  +0.24%  '

$ws.Range("D27").Value = '2.14'
$ws.Range("E27").Value = '  +2.15%  '

$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("D29").Value = '9.99'
$ws.Range("E29").Value = '  +0.28%  '

$ws.Range("D30").Value = '3.913.67'
$ws.Range("E30").Value = '  -0.62%  '

$ws.Range("E31").Value = '  +0.56%  '

$ws.Range("E32").Value = '  +1.76%  '

$ws.Range("D33").Value = '7.26'
$ws.Range("E33").Value = '  -0.13%  '

$ws.Range("D34").Value = '28.79'
$ws.Range("E34").Value = '  -1.11%  '

$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.08%  '

$ws.Range("D36").Value = '3.717.84'
$ws.Range("E36").Value = '  -0.53%  '

$ws.Range("E37").Value = '  -0.31%  '

$ws.Range("D38").Value = '0.157'
$ws.Range("E38").Value = '  +11.97%  '

$ws.Range("E39").Value = '  +2.12%  '

$ws.Range("D40").Value = '3.35'
$ws.Range("E40").Value = '  +2.80%  '

$ws.Range("D41").Value = '5.76'
$ws.Range("E41").Value = '  -0.81%  '

$ws.Range("D42").Value = '0.962'
$ws.Range("E42").Value = '  -2.04%  '

$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.04%  '

$ws.Range("E44").Value = '  +0.01%  '

$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").Value = '154.05'
$ws.Range("E45").Value = '  +1.18%  '

$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").Value = '0.298'
$ws.Range("E46").Value = '  +0.51%  '

$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").Value = '46.62'
$ws.Range("E47").Value = '  -1.56%  '

$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").Value = '42.63'

$ws.Range("E49").Value = '  +2.17%  '

$ws.Range("E50").Value = '  +0.85%  '

$ws.Range("D51").Value = '1.36'
$ws.Range("E51").Value = '  +0.51%  '

# Drop the temporary Text number format again so the touched cells keep the same
# (default) style they had before the edit.
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
